$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "60.202.83"
Set-TextCell "E2" "  +0.08%  "
Set-TextCell "D3" "2.419.16"
Set-TextCell "E3" "  -0.22%  "
Set-TextCell "E4" "  +0.00%  "
Set-TextCell "D5" "554.29"
Set-TextCell "E5" "  -0.06%  "
Set-TextCell "D6" "137.09"
Set-TextCell "E6" "  -1.19%  "
Set-TextCell "D7" "0.999"
Set-TextCell "E8" "  +1.69%  "
Set-TextCell "E9" "  -1.39%  "
Set-TextCell "E10" "  -1.71%  "
Set-TextCell "E11" "  -0.10%  "
Set-TextCell "D12" "0.352"
Set-TextCell "E12" "  -1.80%  "
Set-TextCell "D13" "24.89"
Set-TextCell "E13" "  -0.23%  "
Set-TextCell "D14" "2.849.80"
Set-TextCell "E14" "  -0.17%  "
Set-TextCell "D15" "60.091.88"
Set-TextCell "E15" "  +0.04%  "
Set-TextCell "E16" "  -1.18%  "
Set-TextCell "D17" "2.414.20"
Set-TextCell "E17" "  -0.29%  "
Set-TextCell "D18" "11.26"
Set-TextCell "E18" "  -1.17%  "
Set-TextCell "E19" "  +2.20%  "
Set-TextCell "D20" "327.78"
Set-TextCell "E20" "  -1.50%  "
Set-TextCell "D21" "6.75"
Set-TextCell "E21" "  -0.59%  "
Set-TextCell "E22" "  +0.02%  "
Set-TextCell "D23" "65.36"
Set-TextCell "E23" "  +0.28%  "
Set-TextCell "D24" "0.178"
Set-TextCell "E24" "  +4.00%  "
Set-TextCell "D25" "8.72"
Set-TextCell "E25" "  +1.37%  "
Set-TextCell "D26" "0.999"
Set-TextCell "E26" "  -0.19%  "
Set-TextCell "E27" "  +2.74%  "
Set-TextCell "D28" "0.0₃0772"
Set-TextCell "E28" "  -2.22%  "
Set-TextCell "E29" "  -1.13%  "
Set-TextCell "D30" "170.36"
Set-TextCell "E30" "  +0.93%  "
Set-TextCell "D31" "6.10"
Set-TextCell "E31" "  -3.20%  "
Set-TextCell "D32" "1.07"
Set-TextCell "E32" "  +1.18%  "
Set-TextCell "E33" "  -4.21%  "
Set-TextCell "E35" "  -0.01%  "
Set-TextCell "E36" "  +1.51%  "
Set-TextCell "E37" "  +0.04%  "
Set-TextCell "E38" "  -0.25%  "
Set-TextCell "D39" "327.88"
Set-TextCell "E39" "  +1.95%  "
Set-TextCell "E40" "  -1.09%  "
Set-TextCell "D41" "144.93"
Set-TextCell "E41" "  +3.21%  "
Set-TextCell "E42" "  -1.22%  "
Set-TextCell "B43" "Stellar"
Set-TextCell "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D43" "0.0967"
Set-TextCell "E43" "  +0.56%  "
Set-TextCell "B44" "InjectiveProtocol"
Set-TextCell "C44" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D44" "19.97"
Set-TextCell "E44" "  +2.02%  "
Set-TextCell "D45" "0.0516"
Set-TextCell "E45" "  -1.14%  "
Set-TextCell "E46" "  -0.03%  "
Set-TextCell "D47" "0.0223"
Set-TextCell "E47" "  -1.45%  "
Set-TextCell "E48" "  -0.12%  "
Set-TextCell "E49" "  -3.03%  "
